# Append the September 3rd, 2020 row of raw/clean SSA data to the
# historical log worksheet (row 96), mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 96

# Column A holds dates stored as plain text (shared strings), e.g. "2020-09-02".
# Assigning a date-shaped string directly would make Excel auto-convert it to
# a real date serial value, so we enter it with a leading quote (text-literal
# prefix) to force text, then reset the cell style back to Normal so no
# stray number-format style is left attached to the cell.
$ws.Cells.Item($row, 1).Value = "'2020-09-03"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 616894
$ws.Cells.Item($row, 3).Value = 692319
$ws.Cells.Item($row, 4).Value = 83820
$ws.Cells.Item($row, 5).Value = 66329
$ws.Cells.Item($row, 6).Value = 25.25
